$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range('D2')
$c.NumberFormat = "@"
$c.Value = '34.631.31'
$c.Style = "Normal"
$ws.Range('E2').Value = '  +1.23%  '
$c = $ws.Range('D3')
$c.NumberFormat = "@"
$c.Value = '1.801.41'
$c.Style = "Normal"
$ws.Range('E3').Value = '  +0.96%  '
$ws.Range('E4').Value = '  -0.16%  '
$c = $ws.Range('D5')
$c.NumberFormat = "@"
$c.Value = '227.62'
$c.Style = "Normal"
$ws.Range('E7').Value = '  -0.10%  '
$c = $ws.Range('D8')
$c.NumberFormat = "@"
$c.Value = '32.83'
$c.Style = "Normal"
$ws.Range('E8').Value = '  +2.78%  '
$ws.Range('E9').Value = '  +1.80%  '
$c = $ws.Range('D10')
$c.NumberFormat = "@"
$c.Value = '0.0699'
$c.Style = "Normal"
$ws.Range('E10').Value = '  +1.14%  '
$ws.Range('E11').Value = '  +0.23%  '
$c = $ws.Range('D12')
$c.NumberFormat = "@"
$c.Value = '2.061.50'
$c.Style = "Normal"
$ws.Range('E12').Value = '  +1.04%  '
$c = $ws.Range('D13')
$c.NumberFormat = "@"
$c.Value = '11.21'
$c.Style = "Normal"
$ws.Range('E13').Value = '  +1.35%  '
$c = $ws.Range('D14')
$c.NumberFormat = "@"
$c.Value = '1.796.52'
$c.Style = "Normal"
$ws.Range('E14').Value = '  +0.76%  '
$ws.Range('E15').Value = '  +2.44%  '
$c = $ws.Range('D16')
$c.NumberFormat = "@"
$c.Value = '34.605.47'
$c.Style = "Normal"
$ws.Range('E16').Value = '  +1.24%  '
$ws.Range('E17').Value = '  +2.91%  '
$c = $ws.Range('D18')
$c.NumberFormat = "@"
$c.Value = '68.98'
$c.Style = "Normal"
$ws.Range('E18').Value = '  +1.56%  '
$c = $ws.Range('D19')
$c.NumberFormat = "@"
$c.Value = '0.0₃0807'
$c.Style = "Normal"
$ws.Range('E19').Value = '  +0.89%  '
$c = $ws.Range('D20')
$c.NumberFormat = "@"
$c.Value = '247.76'
$c.Style = "Normal"
$ws.Range('E20').Value = '  +0.23%  '
$ws.Range('E21').Value = '  +2.86%  '
$ws.Range('E22').Value = '  -0.18%  '
$ws.Range('E23').Value = '  +2.27%  '
$c = $ws.Range('D24')
$c.NumberFormat = "@"
$c.Value = '168.02'
$c.Style = "Normal"
$ws.Range('E24').Value = '  +3.54%  '
$ws.Range('E25').Value = '  +1.43%  '
$ws.Range('E26').Value = '  +1.56%  '
$ws.Range('E27').Value = '  +1.86%  '
$ws.Range('E28').Value = '  +2.27%  '
$ws.Range('E29').Value = '  -0.12%  '
$c = $ws.Range('D30')
$c.NumberFormat = "@"
$c.Value = '4.12'
$c.Style = "Normal"
$ws.Range('E30').Value = '  +11.42%  '
$ws.Range('E31').Value = '  +0.95%  '
$ws.Range('B32').Value = 'Filecoin'
$ws.Range('C32').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$c = $ws.Range('D32')
$c.NumberFormat = "@"
$c.Value = '3.82'
$c.Style = "Normal"
$ws.Range('E32').Value = '  +1.90%  '
$ws.Range('B33').Value = 'Hedera'
$ws.Range('C33').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$c = $ws.Range('D33')
$c.NumberFormat = "@"
$c.Value = '0.0525'
$c.Style = "Normal"
$ws.Range('E33').Value = '  +0.88%  '
$ws.Range('E34').Value = '  +2.55%  '
$c = $ws.Range('D35')
$c.NumberFormat = "@"
$c.Value = '1.433.55'
$c.Style = "Normal"
$ws.Range('E35').Value = '  -0.80%  '
$ws.Range('E36').Value = '  +7.92%  '
$c = $ws.Range('D37')
$c.NumberFormat = "@"
$c.Value = '0.674'
$c.Style = "Normal"
$ws.Range('E37').Value = '  +3.13%  '
$ws.Range('E38').Value = '  +3.10%  '
$c = $ws.Range('D39')
$c.NumberFormat = "@"
$c.Value = '0.0192'
$c.Style = "Normal"
$ws.Range('E39').Value = '  +0.40%  '
$c = $ws.Range('D40')
$c.NumberFormat = "@"
$c.Value = '85.41'
$c.Style = "Normal"
$ws.Range('E40').Value = '  +6.28%  '
$ws.Range('E41').Value = '  +1.52%  '
$c = $ws.Range('D42')
$c.NumberFormat = "@"
$c.Value = '0.943'
$c.Style = "Normal"
$ws.Range('E42').Value = '  +2.00%  '
$c = $ws.Range('D43')
$c.NumberFormat = "@"
$c.Value = '2.76'
$c.Style = "Normal"
$ws.Range('E43').Value = '  +3.30%  '
$c = $ws.Range('D44')
$c.NumberFormat = "@"
$c.Value = '13.79'
$c.Style = "Normal"
$ws.Range('E44').Value = '  +0.88%  '
$ws.Range('E45').Value = '  +3.48%  '
$ws.Range('E46').Value = '  +0.54%  '
$ws.Range('E47').Value = '  +0.23%  '
$c = $ws.Range('D48')
$c.NumberFormat = "@"
$c.Value = '1.961.19'
$c.Style = "Normal"
$ws.Range('E48').Value = '  +1.00%  '
$ws.Range('E49').Value = '  +1.46%  '
$ws.Range('E50').Value = '  -0.12%  '
$ws.Range('E51').Value = '  -5.09%  '
